$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 100
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# Hunk 1: ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 986
$ws.Range("I12").Value = 903.2
$ws.Range("K12").Value = 903.2
$ws.Range("M12").Value = -733.2

# Hunk 2: ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3556.6
$ws.Range("I98").Value = 2719.1538
$ws.Range("K98").Value = 2719.1538
$ws.Range("M98").Value = -1221.1538

# Hunk 3: ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1500
$ws.Range("J121").Value = 1500
$ws.Range("L121").Value = 4500
$ws.Range("N121").Value = -7994

# Hunk 4: ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3556.6
$ws.Range("I122").Value = 2719.1538
$ws.Range("K122").Value = 8157.4614
$ws.Range("M122").Value = -5707.4614

# Hunk 5: ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 891.53845
$ws.Range("I135").Value = 780.9091
$ws.Range("K135").Value = 7028.1819
$ws.Range("M135").Value = -4493.1819

# Hunk 6: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4004.9019
$ws.Range("I32").Value = 2608.923
$ws.Range("K32").Value = 2608.923
$ws.Range("M32").Value = -2321.923

# Hunk 7: ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2173.4644
$ws.Range("I61").Value = 1340.48
$ws.Range("K61").Value = 1340.48
$ws.Range("M61").Value = -1128.48

# Hunk 8: ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3209.0715
$ws.Range("J88").Value = 3887.375
$ws.Range("L88").Value = 3887.375
$ws.Range("N88").Value = -4699.375

# Hunk 9: ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3209.0715
$ws.Range("J91").Value = 3887.375
$ws.Range("L91").Value = 3887.375
$ws.Range("N91").Value = -6695.375

# Hunk 10: ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1794.4
$ws.Range("I132").Value = 1565.5
$ws.Range("K132").Value = 4696.5
$ws.Range("M132").Value = -2166.5

# Hunk 11: ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2173.4644
$ws.Range("I136").Value = 1340.48
$ws.Range("K136").Value = 4021.44
$ws.Range("M136").Value = -1471.44

# Hunk 12: BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1325.9
$ws.Range("I20").Value = 1385.625
$ws.Range("J20").Value = 1087
$ws.Range("K20").Value = 1385.625
$ws.Range("L20").Value = 1087
$ws.Range("M20").Value = -1138.625
$ws.Range("N20").Value = -1581

# Hunk 13: BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 134915.27
$ws.Range("I86").Value = 1699.1428
$ws.Range("J86").Value = 251479.38
$ws.Range("K86").Value = 1699.1428
$ws.Range("L86").Value = 251479.38
$ws.Range("M86").Value = -576.1428000000001
$ws.Range("N86").Value = -253725.38

# Hunk 14: BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 134915.27
$ws.Range("I89").Value = 1699.1428
$ws.Range("J89").Value = 251479.38
$ws.Range("K89").Value = 8495.714
$ws.Range("L89").Value = 1257396.9
$ws.Range("M89").Value = -2879.714
$ws.Range("N89").Value = -1268628.9

# Hunk 15: BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 762
$ws.Range("I99").Value = 762
$ws.Range("K99").Value = 762
$ws.Range("M99").Value = 736

# Hunk 16: BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7659.8486
$ws.Range("I134").Value = 9405.429
$ws.Range("J134").Value = 4605.0835
$ws.Range("K134").Value = 28216.287
$ws.Range("L134").Value = 13815.2505
$ws.Range("M134").Value = -25681.287
$ws.Range("N134").Value = -18885.2505

# Hunk 17: CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 821.2857
$ws.Range("I22").Value = 362.25
$ws.Range("J22").Value = 1433.3334
$ws.Range("K22").Value = 362.25
$ws.Range("L22").Value = 1433.3334
$ws.Range("M22").Value = -12.25
$ws.Range("N22").Value = -2133.3334

# Hunk 18: CUL row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 226
$ws.Range("I40").Value = 89
$ws.Range("J40").Value = 500
$ws.Range("K40").Value = 356
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -287
$ws.Range("N40").Value = -2138

# Hunk 19: CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 515.26666
$ws.Range("I107").Value = 281.5
$ws.Range("J107").Value = 551.2308
$ws.Range("K107").Value = 844.5
$ws.Range("L107").Value = 1653.6924
$ws.Range("M107").Value = 1075.5
$ws.Range("N107").Value = -5493.6924

# Hunk 20: CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6033284.5
$ws.Range("I131").Value = 125000570
$ws.Range("J131").Value = 9624.671
$ws.Range("K131").Value = 375001710
$ws.Range("L131").Value = 28874.013
$ws.Range("M131").Value = -374996670
$ws.Range("N131").Value = -38954.013

# Hunk 21: GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5999.6
$ws.Range("I70").Value = 6374.5
$ws.Range("K70").Value = 6374.5
$ws.Range("M70").Value = -6104.5

# Hunk 22: GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5999.6
$ws.Range("I73").Value = 6374.5
$ws.Range("K73").Value = 6374.5
$ws.Range("M73").Value = -5438.5

# Hunk 23: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1994.4814
$ws.Range("I122").Value = 1747.9445
$ws.Range("J122").Value = 2487.5557
$ws.Range("K122").Value = 5243.833500000001
$ws.Range("L122").Value = 7462.6671
$ws.Range("M122").Value = -2793.833500000001
$ws.Range("N122").Value = -12362.6671

# Hunk 24: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4108.2
$ws.Range("I132").Value = 3302.1738
$ws.Range("K132").Value = 9906.5214
$ws.Range("M132").Value = -7376.5214

# Hunk 25: LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8321.429
$ws.Range("I122").Value = 5123.75
$ws.Range("J122").Value = 9600.5
$ws.Range("K122").Value = 15371.25
$ws.Range("L122").Value = 28801.5
$ws.Range("M122").Value = -12921.25
$ws.Range("N122").Value = -33701.5

# Hunk 26: LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2345.6875
$ws.Range("I132").Value = 2266
$ws.Range("J132").Value = 2364.077
$ws.Range("K132").Value = 6798
$ws.Range("L132").Value = 7092.231000000001
$ws.Range("M132").Value = -4268
$ws.Range("N132").Value = -12152.231

# Hunk 27: WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1712.5
$ws.Range("I81").Value = 1712.5
$ws.Range("K81").Value = 3425
$ws.Range("M81").Value = -2364

# Hunk 28: WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1712.5
$ws.Range("I84").Value = 1712.5
$ws.Range("K84").Value = 17125
$ws.Range("M84").Value = -11821

# Hunk 29: WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5516.407
$ws.Range("I126").Value = 4807.684
$ws.Range("K126").Value = 14423.052
$ws.Range("M126").Value = -11953.052

# Hunk 30: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3333.6191
$ws.Range("I132").Value = 3193
$ws.Range("K132").Value = 9579
$ws.Range("M132").Value = -7049

# Hunk 31: WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4348.5
$ws.Range("I136").Value = 4627.7
$ws.Range("J136").Value = 3999.5
$ws.Range("K136").Value = 13883.1
$ws.Range("L136").Value = 11998.5
$ws.Range("M136").Value = -11333.1
$ws.Range("N136").Value = -17098.5
